$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 661, shifting existing rows 661:702 down to 662:703.
$ws.Rows.Item(661).Insert()

# Column A holds date-like text (e.g. "2026/01/18"); momentarily force a
# text number format so the COM layer doesn't auto-coerce the literal into
# a date serial number, then restore the default "Normal" style so the new
# row matches the plain (unstyled) look of its neighbours.
$ws.Cells.Item(661, 1).NumberFormat = "@"
$ws.Cells.Item(661, 1).Value = "2026/01/18"
$ws.Cells.Item(661, 1).Style = "Normal"
$ws.Cells.Item(661, 2).Value = "日"
$ws.Cells.Item(661, 3).Value = 7
$ws.Cells.Item(661, 4).Value = 201
